$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numbers-as-text (e.g. "302.50",
# "1.645.01"); writing them with a leading quote via .Formula forces
# Excel to keep them as literal text instead of renormalizing them
# into a real number (which would drop trailing zeros, etc.) -- the
# same quote-prefix trick a person uses when typing such a value in
# the Excel UI.
function Set-TextCell($range, [string]$text) {
    $range.Formula = "'" + $text
}

Set-TextCell $ws.Range("D2") '23.500.05'
$ws.Range("E2").Value = '  +0.78%  '

Set-TextCell $ws.Range("D3") '1.645.01'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("E5").Value = '  +0.18%  '

Set-TextCell $ws.Range("D6") '302.50'

Set-TextCell $ws.Range("D7") '0.3838'
$ws.Range("E7").Value = '  +0.93%  '

$ws.Range("E8").Value = '  +1.06%  '

Set-TextCell $ws.Range("D9") '50.98'
$ws.Range("E9").Value = '  -1.88%  '

Set-TextCell $ws.Range("D10") '0.08159'
$ws.Range("E10").Value = '  +0.44%  '

$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("E12").Value = '  +0.14%  '

Set-TextCell $ws.Range("D13") '22.26'
$ws.Range("E13").Value = '  +0.26%  '

Set-TextCell $ws.Range("D14") '6.427'
$ws.Range("E14").Value = '  +0.28%  '

Set-TextCell $ws.Range("D15") '7.422'
$ws.Range("E15").Value = '  +1.99%  '

Set-TextCell $ws.Range("D16") '0.00001216'
$ws.Range("E16").Value = '  -0.58%  '

Set-TextCell $ws.Range("D17") '1.650.49'
$ws.Range("E17").Value = '  +1.40%  '

Set-TextCell $ws.Range("D18") '97.30'
$ws.Range("E18").Value = '  +2.65%  '

Set-TextCell $ws.Range("D19") '0.07011'
$ws.Range("E19").Value = '  +1.09%  '

Set-TextCell $ws.Range("D20") '6.752'
$ws.Range("E20").Value = '  +3.07%  '

$ws.Range("E21").Value = '  +1.27%  '

Set-TextCell $ws.Range("D22") '1.004'
$ws.Range("E22").Value = '  +0.18%  '

Set-TextCell $ws.Range("D23") '12.60'
$ws.Range("E23").Value = '  +1.73%  '

Set-TextCell $ws.Range("D24") '23.502.07'
$ws.Range("E24").Value = '  +0.74%  '

Set-TextCell $ws.Range("D25") '2.483'
$ws.Range("E25").Value = '  -1.98%  '

Set-TextCell $ws.Range("D26") '3.035'
$ws.Range("E26").Value = '  -2.34%  '

Set-TextCell $ws.Range("D27") '21.17'
$ws.Range("E27").Value = '  +0.48%  '

Set-TextCell $ws.Range("D28") '153.69'
$ws.Range("E28").Value = '  +1.45%  '

Set-TextCell $ws.Range("D29") '5.227'
$ws.Range("E29").Value = '  -0.60%  '

Set-TextCell $ws.Range("D30") '133.83'
$ws.Range("E30").Value = '  +0.83%  '

Set-TextCell $ws.Range("D31") '1.831.38'
$ws.Range("E31").Value = '  +1.15%  '

Set-TextCell $ws.Range("D32") '7.040'
$ws.Range("E32").Value = '  +8.58%  '

Set-TextCell $ws.Range("D33") '2.254'
$ws.Range("E33").Value = '  +5.17%  '

Set-TextCell $ws.Range("D34") '12.19'
$ws.Range("E34").Value = '  +6.62%  '

Set-TextCell $ws.Range("D35") '1.053'

Set-TextCell $ws.Range("D36") '0.02784'
$ws.Range("E36").Value = '  +1.80%  '

$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range("D37") '0.2492'
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range("D38") '0.08780'
$ws.Range("E38").Value = '  +0.62%  '

Set-TextCell $ws.Range("D39") '6.053'
$ws.Range("E39").Value = '  +2.17%  '

Set-TextCell $ws.Range("D40") '0.06964'
$ws.Range("E40").Value = '  +0.13%  '

Set-TextCell $ws.Range("D41") '13.09'
$ws.Range("E41").Value = '  +8.14%  '

Set-TextCell $ws.Range("D42") '0.6961'
$ws.Range("E42").Value = '  +0.26%  '

Set-TextCell $ws.Range("D43") '1.332'
$ws.Range("E43").Value = '  +0.55%  '

Set-TextCell $ws.Range("D44") '15.93'
$ws.Range("E44").Value = '  +3.02%  '

Set-TextCell $ws.Range("D45") '0.6475'
$ws.Range("E45").Value = '  +1.21%  '

Set-TextCell $ws.Range("D46") '1.003'
$ws.Range("E46").Value = '  +0.19%  '

Set-TextCell $ws.Range("D47") '2.285'
$ws.Range("E47").Value = '  +1.13%  '

Set-TextCell $ws.Range("D49") '0.07867'
$ws.Range("E49").Value = '  -0.65%  '

Set-TextCell $ws.Range("D50") '127.74'
$ws.Range("E50").Value = '  -1.40%  '

$ws.Range("E51").Value = '  -0.27%  '
